$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Append a new row (row 60) with the next quarterly period and its values.
# Force text format on A60 so the date-like string "01-07-2021" is stored
# as a plain text string (shared string), matching the rest of column A,
# instead of being auto-converted by Excel into a date serial number.
$ws.Range("A60").NumberFormat = "@"
$ws.Range("A60").Value = "01-07-2021"
$ws.Range("A60").Style = "Normal"
$ws.Range("B60").Value = 9813
$ws.Range("C60").Value = 2482
$ws.Range("D60").Value = 7331
$ws.Range("E60").Value = 7353
$ws.Range("F60").Value = 3307
$ws.Range("G60").Value = 4046
$ws.Range("H60").Value = 202
$ws.Range("I60").Value = 202
$ws.Range("J60").Value = 584
$ws.Range("K60").Value = 257
